$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update raw values for rows 4-8 per the "Legs and Room 3 Sesi 2" update
$ws.Range("C4").Value = 1500
$ws.Range("D4").Value = 1550
$ws.Range("F4").Value = 2000
$ws.Range("G4").Value = 2100

$ws.Range("D5").Value = 1450

$ws.Range("C6").Value = 1350
$ws.Range("D6").Value = 1650

$ws.Range("D7").Value = 1350
$ws.Range("G7").Value = 850

$ws.Range("C8").Value = 1550
$ws.Range("D8").Value = 1300
$ws.Range("E8").Value = 1870

# Update view state: scroll to show column E onward, and move selection to F14
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("F14").Select()
